# "Corrected J assignment (typo in ASD)"
#
# The worksheet pulls observed term-energy differences (column E, via the
# n=10..12 block in rows 18-21) into the per-series columns I:L in row 9.
# The J=(1/2,1/2);0 series (column K) and J=(1/2,3/2);1 series (column L)
# had been cross-wired: K9 pointed at row 20 (n=11) and L9 at row 21
# (n=10) instead of the other way around. This also means the "n" index
# typed into C20/C21 was swapped, and the value that should have carried
# the "typo in ASD" flag belongs on row 21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Swap the K9/L9 formulas so they pick up the correct source rows.
$ws.Range("K9").Formula = '=E21'
$ws.Range("L9").Formula = '=E20'

# 2. Swap the C20/C21 "n" labels (and the associated red-font highlight
#    style) back to the correct order.
$ws.Range("C20").ClearFormats()
$ws.Range("C21").Value = 1
$ws.Range("C20").Value = 2
$ws.Range("C21").Font.Color = 255

# 3. Leave a note next to the corrected row explaining the fix, using the
#    same red-font style used to flag the correction.
$ws.Range("F21").Value = "typo in ASD"
$ws.Range("F21").Font.Color = 255

# 4. Refresh the dependent formula ranges so everything recalculates from
#    the corrected data (re-entering them also re-normalises the shared
#    formula groupings the way Excel does on edit).
$ws.Range("E6:E31").Formula = '=IE-D6'
$ws.Range("N6:Q11").Formula = '=Rhc/($H6-N$4)^2'
$ws.Range("S7:S11").Formula = '=N7-I7'
$ws.Range("T6:T11").Formula = '=O6-J6'
$ws.Range("U6:U12").Formula = '=P6-K6'
$ws.Range("V6:V12").Formula = '=Q6-L6'
$ws.Range("O13:Q13").Formula = '=SQRT(SUMXMY2(O6:O12,J6:J12)/COUNT(O6:O12))'
$ws.Range("I16:L21").Formula = '=$H16-SQRT(Rhc/I6)'

# 5. Re-solve for the quantum defects (P4, Q4) now that the underlying
#    data has been corrected (previously minimised via Solver).
$ws.Range("P4").Value = 4.3017486358129604
$ws.Range("Q4").Value = 4.2918096623756661

# 6. Leave the selection where the fix was made.
$ws.Range("U9").Select()
